# Updated symbol list on Fri Jan 27 15:38:54 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) values on the
# cryptos sheet. Values are stored as text in the workbook (e.g. "304.53",
# "-0.97%"), so we temporarily mark each cell as text ("@") before writing
# the value -- otherwise Excel would auto-convert these strings into
# numbers/percentages. We then restore NumberFormat/Style so the cell's
# formatting metadata stays identical to the original (no stray styles).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "304.53"
Set-TextValue $ws.Range("E2") "-0.97%"

Set-TextValue $ws.Range("D3") "35.79"
Set-TextValue $ws.Range("E3") "-0.42%"

Set-TextValue $ws.Range("D4") "5.044"
Set-TextValue $ws.Range("E4") "-0.49%"

Set-TextValue $ws.Range("D5") "0.08041"
Set-TextValue $ws.Range("E5") "-0.86%"

Set-TextValue $ws.Range("D6") "1.859"
Set-TextValue $ws.Range("E6") "-4.14%"

Set-TextValue $ws.Range("D7") "4.119"
Set-TextValue $ws.Range("E7") "-1.10%"

Set-TextValue $ws.Range("D8") "7.774"
Set-TextValue $ws.Range("E8") "-0.64%"

Set-TextValue $ws.Range("D9") "0.9259"
Set-TextValue $ws.Range("E9") "-1.34%"

Set-TextValue $ws.Range("E10") "-7.44%"

Set-TextValue $ws.Range("D11") "0.1887"
Set-TextValue $ws.Range("E11") "-1.23%"

Set-TextValue $ws.Range("D12") "0.09064"
Set-TextValue $ws.Range("E12") "-1.72%"

Set-TextValue $ws.Range("D13") "0.03440"
Set-TextValue $ws.Range("E13") "-1.69%"

Set-TextValue $ws.Range("D14") "0.09868"
Set-TextValue $ws.Range("E14") "-0.39%"

Set-TextValue $ws.Range("D15") "0.001413"
Set-TextValue $ws.Range("E15") "-2.46%"

Set-TextValue $ws.Range("D16") "0.006266"
Set-TextValue $ws.Range("E16") "7.06%"

Set-TextValue $ws.Range("D17") "3.863"
Set-TextValue $ws.Range("E17") "6.55%"

Set-TextValue $ws.Range("D19") "0.3409"
Set-TextValue $ws.Range("E19") "-0.58%"

Set-TextValue $ws.Range("D20") "0.1299"
Set-TextValue $ws.Range("E20") "-3.45%"

Set-TextValue $ws.Range("D21") "4.822"
Set-TextValue $ws.Range("E21") "-7.07%"

Set-TextValue $ws.Range("D22") "0.2406"
Set-TextValue $ws.Range("E22") "-5.02%"

Set-TextValue $ws.Range("D23") "0.04367"
Set-TextValue $ws.Range("E23") "-0.82%"

Set-TextValue $ws.Range("D24") "0.001230"
Set-TextValue $ws.Range("E24") "-0.32%"

Set-TextValue $ws.Range("D25") "0.004841"
Set-TextValue $ws.Range("E25") "1.51%"

Set-TextValue $ws.Range("D27") "0.0001300"
Set-TextValue $ws.Range("E27") "0.04%"

Set-TextValue $ws.Range("E28") "42.15%"

Set-TextValue $ws.Range("D39") "0.01970"
Set-TextValue $ws.Range("E39") "-3.16%"

Set-TextValue $ws.Range("D40") "0.05142"
Set-TextValue $ws.Range("E40") "0.65%"

Set-TextValue $ws.Range("D41") "0.007522"
Set-TextValue $ws.Range("E41") "-1.48%"

Set-TextValue $ws.Range("D42") "0.01008"
Set-TextValue $ws.Range("E42") "-10.43%"

Set-TextValue $ws.Range("D43") "0.1354"
Set-TextValue $ws.Range("E43") "-1.88%"

Set-TextValue $ws.Range("D44") "0.002111"
Set-TextValue $ws.Range("E44") "0.52%"

Set-TextValue $ws.Range("D45") "0.009869"
Set-TextValue $ws.Range("E45") "-12.89%"

Set-TextValue $ws.Range("D46") "0.00006187"
Set-TextValue $ws.Range("E46") "-2.21%"

Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "0.00%"

Set-TextValue $ws.Range("D48") "64.85"
Set-TextValue $ws.Range("E48") "-0.17%"

Set-TextValue $ws.Range("E49") "4.84%"

Set-TextValue $ws.Range("D50") "0.00002100"
Set-TextValue $ws.Range("E50") "0.00%"

Set-TextValue $ws.Range("D51") "0.0002000"
Set-TextValue $ws.Range("E51") "0.00%"
